$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map each emoji "statut" marker to its plain-text replacement
# (the author's fix for an Excel/emoji rendering problem).
$map = @{
    "📕" = "-3";
    "📘" = "⚠️";
    "📙" = "+3";
    "📗" = "✅"
}

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

# Column A ("statut") holds the emoji markers. Rewrite every matching
# cell as a quoted-text formula first (so values that look numeric,
# like "-3"/"+3", are stored as text rather than being coerced into
# numbers), then flatten the whole column back down to plain values.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value2
    if ($map.ContainsKey($old)) {
        $new = $map[$old]
        $cell.Formula = '="' + $new + '"'
    }
}

$rng = $ws.Range("A2:A" + $lastRow)
$rng.Copy()
$rng.PasteSpecial(-4163)
$excel.CutCopyMode = $false
